$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows 2-8 (columns B, C, E, F, G change; A and D stay the same)
# Note: numbers with magnitude requiring scientific notation are written out in plain
# decimal form (still the exact same double value) because the script parser does not
# accept "E" exponent notation in numeric literals.
$data = @(
    @(2, 1,    0.009122927242422047,  0.83910583052295618,      0.001183986663818359, 0.99236526923462420),
    @(3, 11,   0.05130013785025216,   0.09525825931251762,      0.004054069519042969, 0.96716117625147124),
    @(4, 19,   0.05643587430324096,   0.004683158876088452,     0.006312370300292969, 0.93901886456140216),
    @(5, 95,   0.05672714649778521,   0.0004537885737423669,    0.02941203117370605,  0.85132398763639561),
    @(6, 115,  0.05669951160987683,   0.00003358704818894733,   0.03550863265991211,  0.83833013848788551),
    @(7, 312,  0.05670125280906785,   0.000002878834664001203,  0.09159636497497559,  0.75344856228451795),
    @(8, 594,  0.05670139881194544,   0.0000003038923118988512, 0.1734015941619873,   0.68448557916741914)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

# Add new row 9 with data, matching style of other "A" column cells (bold/centered/bordered)
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(8, 1).Copy()
$ws.Cells.Item(9, 1).PasteSpecial(-4122)

$ws.Cells.Item(9, 2).Value = 3198
$ws.Cells.Item(9, 3).Value = 0.05670141636398267
$ws.Cells.Item(9, 4).Value = 0.0000001
$ws.Cells.Item(9, 5).Value = 0.000000005659696698042492
$ws.Cells.Item(9, 6).Value = 0.950007438659668
$ws.Cells.Item(9, 7).Value = 0.4519505583416418
